# follow_arrival.xlsx: rename the "form_id" setting to "table_id" and add a
# new "properties" sheet describing the column order, so the downstream
# tooling can emit definitions.csv / properties.csv.

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item(1)   # "survey"
$settings = $wb.Worksheets.Item(2)   # "settings"

# --- settings sheet: form_id -> table_id (value "follow_arrival" stays) ---
$settings.Cells.Item(2, 1).Value = "table_id"

# --- add the new "properties" sheet after the last existing sheet ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$properties = $wb.Worksheets.Add($null, $lastSheet)
$properties.Name = "properties"

# header row
$properties.Cells.Item(1, 1).Value = "partition"
$properties.Cells.Item(1, 2).Value = "aspect"
$properties.Cells.Item(1, 3).Value = "key"
$properties.Cells.Item(1, 4).Value = "type"
$properties.Cells.Item(1, 5).Value = "value"

# data row: the table's column order, stored as a JSON array
$properties.Cells.Item(2, 1).Value = "Table"
$properties.Cells.Item(2, 2).Value = "default"
$properties.Cells.Item(2, 3).Value = "colOrder"
$properties.Cells.Item(2, 4).Value = "array"
$properties.Cells.Item(2, 5).Value = '["FA_FOL_date","FA_FOL_B_focal_AnimID","FA_B_arr_AnimID","FA_seq_num","FA_type_of_certainty","FA_type_of_nesting","FA_type_of_cycle","FA_time_start","FA_time_end","FA_duration_of_obs","FA_within_five_meters","FA_closest_to_focal"]'

# --- view/selection state: settings keeps A3 selected, properties becomes
#     the active tab with E7 selected (survey's C14 selection is untouched) ---
$settings.Activate() | Out-Null
$settings.Range("A3").Select() | Out-Null

$properties.Activate() | Out-Null
$properties.Range("E7").Select() | Out-Null
